$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column is inserted before the old "POI-keys" column (old F), pushing the
# existing POI-keys data (and everything to its right) one column over.
$ws.Columns("F").Insert()

# Row 1 - headers
$ws.Range("E1").Value = "Wegzweck"
$ws.Range("F1").Value = "km_Wegzweck"

# Row 2 - Summe Schule
$ws.Range("E2").Value = "Ausbildung;Begleitung"
$ws.Range("F2").Formula = "=(7+8)/2"

# Row 3 - Summe Einkauf
$ws.Range("E3").Value = "Einkauf"
$ws.Range("F3").Value = 5

# Row 4 - Eigener Arbeitsplatz
$ws.Range("E4").Value = "Arbeit"
$ws.Range("F4").Value = 12

# Row 5 - Anderer Dienstort/-weg
$ws.Range("E5").Value = "Arbeit"
$ws.Range("F5").Value = 18

# Row 6 - Behördengang. Arztbesuch
$ws.Range("E6").Value = "Erledigung"
$ws.Range("F6").Value = 10

# Row 7 - Dienstleistungseinrichtung (z. B. Post. Bank. Friseur)
$ws.Range("E7").Value = "Erledigung"
$ws.Range("F7").Value = 10

# Row 8 - Kultur. Theater. Kino
$ws.Range("E8").Value = "Freizeit"
$ws.Range("F8").Value = 15

# Row 9 - Gaststätte/Kneipe
$ws.Range("E9").Value = "Freizeit"
$ws.Range("F9").Value = 15

# Row 10 - Summe Sport Freizeit
$ws.Range("E10").Value = "Freizeit"
$ws.Range("F10").Value = 15

# Row 11 - Eigene Wohnung
$ws.Range("E11").Value = "Freizeit"
$ws.Range("F11").Value = 15

# Selection moves to E2, matching the edited workbook
$ws.Range("E2").Select() | Out-Null
